# Apply ranking corrections to the davis_2019 match sheet.
# Column V = single_ranking_player1, W = doubles_ranking_player1,
# Y = single_ranking_player2, Z = doubles_ranking_player2.
#
# These columns store rankings as TEXT (e.g. "536=" for a tied ranking),
# so purely-numeric replacement values are entered with a leading
# apostrophe to keep Excel from auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V2").Value = "'631"

$ws.Range("V5").Value = "'22"

$ws.Range("V7").Value = "'22"

$ws.Range("V8").Value = "'951"
$ws.Range("W8").Value = "536="

$ws.Range("V9").Value = "'612"
$ws.Range("W9").Value = "1385="

$ws.Range("W10").Value = "'342"

$ws.Range("W12").Value = "'54"
$ws.Range("Y12").Value = "'32"

$ws.Range("V13").Value = "'612"
$ws.Range("W13").Value = "1385="

$ws.Range("V14").Value = "'452"

$ws.Range("V16").Value = "'184"
$ws.Range("W16").Value = "'868"

$ws.Range("V17").Value = "'35"
$ws.Range("W17").Value = "'131"

$ws.Range("V18").Value = "1079="

$ws.Range("V19").Value = "'213"

$ws.Range("V20").Value = "'244"

$ws.Range("V21").Value = "'154"
$ws.Range("W21").Value = "'863"

$ws.Range("W22").Value = "'33"

$ws.Range("Y23").Value = "'154"
$ws.Range("Z23").Value = "'863"

$ws.Range("V25").Value = "'612"
$ws.Range("W25").Value = "1385="

$ws.Range("V26").Value = "'22"

$ws.Range("V28").Value = "'22"

$ws.Range("V29").Value = "'612"
$ws.Range("W29").Value = "1385="

$ws.Range("V30").Value = "'452"

$ws.Range("V32").Value = "'184"
$ws.Range("W32").Value = "'868"

$ws.Range("V33").Value = "'154"
$ws.Range("W33").Value = "'863"

$ws.Range("W34").Value = "'32"

$ws.Range("Y35").Value = "'154"
$ws.Range("Z35").Value = "'863"

$ws.Range("W36").Value = "'550"

$ws.Range("V39").Value = "'154"
$ws.Range("W39").Value = "'863"

Write-Host "Applied 42 ranking cell updates"
